$d = $word.ActiveDocument

function Replace-ParagraphText($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range.Duplicate
    # Wrap = 0 (wdFindStop) keeps the search confined to this paragraph's
    # range so identical text elsewhere in the document can never be hit.
    $ok = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                             $true, 0, $false, $newText, 2)
    if (-not $ok) {
        throw "Find failed for paragraph $paraIndex expecting '$oldText'"
    }
}

# Title (Heading1) and its bold duplicate later in the document
Replace-ParagraphText 1 "Play Multifly Free - Review of Exciting Jungle Themed Slot" "Play Multifly Free - Exciting Gameplay and Potential Wins"
Replace-ParagraphText 43 "Play Multifly Free - Review of Exciting Jungle Themed Slot" "Play Multifly Free - Exciting Gameplay and Potential Wins"

# "What we like" bullet list (paragraphs 37-39) - values rotate among themselves
Replace-ParagraphText 37 "Chameleon multiplier meters" "Potential win of 10,000 times the bet"
Replace-ParagraphText 38 "Vibrant colors and well-designed symbols" "243 ways to win"
Replace-ParagraphText 39 "Potential win of 10,000x the bet" "High-quality design"

# "What we don't like" bullet list (paragraphs 41-42) - values rotate among themselves
Replace-ParagraphText 41 "May not be suitable for less experienced players" "High volatility"
Replace-ParagraphText 42 "High volatility" "May be better suited for experienced players"

# Italic meta description paragraph
Replace-ParagraphText 44 "Read our review of Multifly, a high-volatility, jungle-themed slot game with chameleon multipliers, dropdown menus, and potential wins of 10,000x the bet. Play for free now." "Play Multifly for free and experience exciting gameplay features and the potential for big wins."

Write-Output "Done"
